$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "64.746.59"
Set-TextValue "E2" "  -0.72%  "

Set-TextValue "D3" "3.446.70"
Set-TextValue "E3" "  -1.13%  "

Set-TextValue "E4" "  -0.02%  "

Set-TextValue "D5" "573.51"
Set-TextValue "E5" "  -1.19%  "

Set-TextValue "D6" "159.17"
Set-TextValue "E6" "  -2.37%  "

Set-TextValue "E7" "  +0.01%  "

Set-TextValue "D8" "3.444.72"
Set-TextValue "E8" "  -1.23%  "

Set-TextValue "D9" "0.575"
Set-TextValue "E9" "  -6.31%  "

Set-TextValue "D10" "7.21"
Set-TextValue "E10" "  -0.76%  "

Set-TextValue "D11" "0.121"
Set-TextValue "E11" "  -3.57%  "

Set-TextValue "D12" "0.438"
Set-TextValue "E12" "  -2.08%  "

Set-TextValue "D13" "4.043.30"
Set-TextValue "E13" "  -1.10%  "

Set-TextValue "E14" "  -0.45%  "

Set-TextValue "D15" "27.55"
Set-TextValue "E15" "  -4.12%  "

Set-TextValue "D16" "0.0000174"
Set-TextValue "E16" "  -10.31%  "

Set-TextValue "D17" "64.786.68"
Set-TextValue "E17" "  -0.70%  "

Set-TextValue "D18" "3.462.09"
Set-TextValue "E18" "  -0.86%  "

Set-TextValue "D19" "6.18"
Set-TextValue "E19" "  -4.39%  "

Set-TextValue "D20" "13.69"
Set-TextValue "E20" "  -4.81%  "

Set-TextValue "D21" "376.83"
Set-TextValue "E21" "  -1.91%  "

Set-TextValue "D22" "7.92"
Set-TextValue "E22" "  -3.77%  "

Set-TextValue "D24" "72.13"
Set-TextValue "E24" "  -0.71%  "

Set-TextValue "D25" "0.533"
Set-TextValue "E25" "  -3.89%  "

Set-TextValue "E26" "  -1.10%  "

Set-TextValue "D27" "9.92"
Set-TextValue "E27" "  -1.66%  "

Set-TextValue "E28" "  -0.26%  "

Set-TextValue "D29" "0.999"
Set-TextValue "E29" "  -0.04%  "

Set-TextValue "E30" "  -6.58%  "

Set-TextValue "D31" "6.06"
Set-TextValue "E31" "  -1.84%  "

Set-TextValue "E32" "  -2.59%  "

Set-TextValue "D33" "23.15"
Set-TextValue "E33" "  -2.49%  "

Set-TextValue "D34" "6.97"
Set-TextValue "E34" "  -3.11%  "

Set-TextValue "E35" "  -3.73%  "

Set-TextValue "D36" "161.16"

Set-TextValue "E37" "  -2.65%  "

Set-TextValue "D38" "2.888.68"
Set-TextValue "E38" "  -3.81%  "

Set-TextValue "D39" "0.0746"
Set-TextValue "E39" "  -4.60%  "

Set-TextValue "D40" "26.16"
Set-TextValue "E40" "  -2.74%  "

Set-TextValue "B41" "OKB"
Set-TextValue "C41" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D41" "43.01"
Set-TextValue "E41" "  +0.07%  "

Set-TextValue "B42" "Mantle"
Set-TextValue "C42" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D42" "0.788"
Set-TextValue "E42" "  +0.90%  "

Set-TextValue "B43" "Filecoin"
Set-TextValue "C43" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D43" "4.51"
Set-TextValue "E43" "  -1.79%  "

Set-TextValue "D44" "26.15"
Set-TextValue "E44" "  +0.71%  "

Set-TextValue "E45" "  -6.23%  "

Set-TextValue "D46" "0.0309"
Set-TextValue "E46" "  -4.20%  "

Set-TextValue "D47" "2.40"
Set-TextValue "E47" "  +9.35%  "

Set-TextValue "D48" "320.79"
Set-TextValue "E48" "  -1.00%  "

Set-TextValue "D49" "1.08"
Set-TextValue "E49" "  -3.32%  "

Set-TextValue "D50" "6.46"
Set-TextValue "E50" "  -4.17%  "

Set-TextValue "D51" "0.844"
Set-TextValue "E51" "  -4.29%  "
